$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 2023 column (K) data
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1220
$ws.Range("K5").Value = 452
$ws.Range("K6").Value = 768

# Copy styles from J column cells to K column cells so formatting matches
$ws.Range("J3").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("J4").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null

$ws.Range("J5").Copy() | Out-Null
$ws.Range("K5").PasteSpecial(-4122) | Out-Null

$ws.Range("J6").Copy() | Out-Null
$ws.Range("K6").PasteSpecial(-4122) | Out-Null
